$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E4").Value2 = -0.04775277527752775
$ws.Range("F4").Value2 = 0.1504793935334936
$ws.Range("E5").Value2 = -0.1068946894689469
$ws.Range("F5").Value2 = -0.1958707782888
$ws.Range("E6").Value2 = -0.069006900690069
$ws.Range("F6").Value2 = 0.0686763731707913
$ws.Range("E7").Value2 = -0.07122712271227123
$ws.Range("F7").Value2 = 0.07346557271634362
$ws.Range("E8").Value2 = -0.1457065706570657
$ws.Range("F8").Value2 = 0.102462679339805
$ws.Range("E9").Value2 = 0.1011821182118212
$ws.Range("F9").Value2 = -0.05912915376422802
$ws.Range("E10").Value2 = 0.1719651965196519
$ws.Range("F10").Value2 = 0.07548601627462351
$ws.Range("E11").Value2 = -0.03737173717371737
$ws.Range("F11").Value2 = -0.02027926682569811
$ws.Range("E12").Value2 = -0.08195619561956194
$ws.Range("F12").Value2 = -0.07302905713276464
$ws.Range("E13").Value2 = -0.004968496849684968
$ws.Range("F13").Value2 = -0.06212240348019821
$ws.Range("E14").Value2 = -0.08415241524152414
$ws.Range("F14").Value2 = 0.05662854306401124
$ws.Range("E15").Value2 = 0.05249324932493248
$ws.Range("F15").Value2 = -0.08649868085463053
$ws.Range("E16").Value2 = -0.1452145214521452
$ws.Range("F16").Value2 = 0.07412658202862039
$ws.Range("E17").Value2 = 0.01056105610561056
$ws.Range("F17").Value2 = 0.13242012024714
$ws.Range("E18").Value2 = -0.09276927692769275
$ws.Range("F18").Value2 = -0.114073993863006
$ws.Range("E19").Value2 = 0.4360876087608761
$ws.Range("F19").Value2 = 0.04221729286899638
$ws.Range("E20").Value2 = 0.0999099909990999
$ws.Range("F20").Value2 = -0.01784725143147232
$ws.Range("E21").Value2 = 0.4471047104710471
$ws.Range("F21").Value2 = 0.006011443179573488
$ws.Range("E22").Value2 = -0.02547854785478548
$ws.Range("F22").Value2 = -0.05240057784020333
$ws.Range("E23").Value2 = 0.4794239423942394
$ws.Range("F23").Value2 = 0.04380122084369728
$ws.Range("E24").Value2 = 0.5956195619561956
$ws.Range("F24").Value2 = 0.0221251041505464
$ws.Range("E25").Value2 = 0.0958175817581758
$ws.Range("F25").Value2 = 0.1731844398165092
$ws.Range("E26").Value2 = 0.113987398739874
$ws.Range("F26").Value2 = -0.03755904799851775
$ws.Range("E27").Value2 = -0.1024902490249025
$ws.Range("F27").Value2 = 0.009285310056415896
$ws.Range("E28").Value2 = 0.0257065706570657
$ws.Range("F28").Value2 = 0.05297428403576429
$ws.Range("E29").Value2 = 0.07840384038403839
$ws.Range("F29").Value2 = -0.2167112794362426
$ws.Range("E30").Value2 = -0.01545754575457545
$ws.Range("F30").Value2 = 0.09877724062701664
$ws.Range("E31").Value2 = -0.03396339633963396
$ws.Range("F31").Value2 = -0.06486621571983757
$ws.Range("E32").Value2 = -0.05942994299429943
$ws.Range("F32").Value2 = -0.06449205950534129
$ws.Range("E33").Value2 = -0.0393039303930393
$ws.Range("F33").Value2 = 0.1006916732578559
$ws.Range("E34").Value2 = -0.2019081908190819
$ws.Range("F34").Value2 = 0.2960136890987281
$ws.Range("E35").Value2 = 0.03368736873687368
$ws.Range("F35").Value2 = 0.00846840232143236
$ws.Range("E36").Value2 = 0.1065226522652265
$ws.Range("F36").Value2 = 0.002575441943116027
$ws.Range("E37").Value2 = 0.07654365436543654
$ws.Range("F37").Value2 = 0.1156267421531656
